$d = $word.ActiveDocument

# --- Change 1: Condense the three CORE COMPETENCIES detail paragraphs into one ---
# The paragraphs hold the "Data Visualization & Design: ...",
# "Geospatial Analysis & Mapping: ...", and "Technical Visualization: ..." text
# right after the CORE COMPETENCIES heading. Delete the last two (including
# their paragraph marks), then overwrite the first paragraph's text with the
# condensed summary line.
$headingRange = $d.Content
$headingFound = $headingRange.Find.Execute("CORE COMPETENCIES", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$pDesign = $d.Paragraphs.Item($headingRange.Paragraphs.First.Index + 1)
$pGeo = $d.Paragraphs.Item($pDesign.Index + 1)
$pTech = $d.Paragraphs.Item($pDesign.Index + 2)

$removeRange = $d.Range($pGeo.Range.Start, $pTech.Range.End)
$removeRange.Delete()

$pDesign.Range.Text = "Data Visualization & Design • Geospatial Analysis & Mapping • Technical Visualization"

# --- Change 2: Add a new TECHNICAL SKILLS section near the end of the document ---
# Locate the "Led multi-million dollar research projects" bullet (the last
# bullet before the closing "For a more detailed..." line) and insert the new
# heading plus three detail paragraphs right after it.
$anchorRange = $d.Content
$found = $anchorRange.Find.Execute("Led multi-million dollar research projects with focus on visual communication of insights and findings", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Create all four new (empty) paragraphs first, anchored right after the
# found bullet, then fill them in by index -- re-using the same Find range
# for InsertParagraphAfter always inserts immediately after the anchor, so
# doing every insert before any text assignment keeps the paragraphs in the
# intended top-to-bottom order.
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()

$anchorIndex = $anchorRange.Paragraphs.First.Index

$pHeading = $d.Paragraphs.Item($anchorIndex + 1)
$pHeading.Range.Text = "TECHNICAL SKILLS"
$pHeading.Range.Style = "Heading2"

$pDataViz = $d.Paragraphs.Item($anchorIndex + 2)
$pDataViz.Range.Text = "DATA VISUALIZATION & DESIGN Interactive Dashboards; Statistical Visualization; Geospatial Mapping; Choropleth Design"

$pGeoSkills = $d.Paragraphs.Item($anchorIndex + 3)
$pGeoSkills.Range.Text = "GEOSPATIAL ANALYSIS & MAPPING Spatial Analysis; Mapping Technologies; Web Mapping; Spatial Data Processing"

$pTechSkills = $d.Paragraphs.Item($anchorIndex + 4)
$pTechSkills.Range.Text = "TECHNICAL VISUALIZATION Programming; Database Integration; Web Technologies; Statistical Computing"

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
